# Updated cryptos list on Mon Sep  4 19:43:09 UTC 2023 with GitHub Actions
#
# Refreshes the Price / Volume(1h) columns with freshly scraped values. Two
# coins (rows 12-13) swapped rank order, so their Coin name + Link also move.
#
# NumberFormat is forced to Text ("@") before writing any value that looks
# numeric (e.g. "215.91", "4.257") because the sheet stores Price/Volume as
# plain text (some prices even contain two "." like "25.884.13") and Excel
# would otherwise silently coerce the assignment into a numeric cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '25.884.13'
$ws.Range("E2").Value = '  -0.45%  '
# Row 3
$ws.Range("D3").Value = '1.631.65'
$ws.Range("E3").Value = '  -0.51%  '
# Row 4
$ws.Range("E4").Value = '  -0.09%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.91'
$ws.Range("E5").Value = '  +0.57%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5116'
$ws.Range("E6").Value = '  +0.49%  '
# Row 7
$ws.Range("E7").Value = '  -0.01%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2565'
$ws.Range("E8").Value = '  -0.04%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06333'
$ws.Range("E9").Value = '  -0.60%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.44'
$ws.Range("E10").Value = '  -0.32%  '
# Row 11
$ws.Range("E11").Value = '  +0.21%  '
# Row 12
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.257'
$ws.Range("E12").Value = '  -0.64%  '
# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.638.03'
$ws.Range("E13").Value = '  -0.93%  '
# Row 14
$ws.Range("D14").Value = '1.855.79'
$ws.Range("E14").Value = '  -0.58%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5503'
$ws.Range("E15").Value = '  +1.06%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.73'
$ws.Range("E16").Value = '  -0.79%  '
# Row 17
$ws.Range("D17").Value = '0.0₅7633'
$ws.Range("E17").Value = '  -1.42%  '
# Row 18
$ws.Range("D18").Value = '25.902.65'
$ws.Range("E18").Value = '  -0.43%  '
# Row 19
$ws.Range("E19").Value = '  +0.02%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '194.50'
$ws.Range("E20").Value = '  -0.95%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.412'
$ws.Range("E21").Value = '  -0.23%  '
# Row 22
$ws.Range("E22").Value = '  -0.75%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.019'
$ws.Range("E23").Value = '  -0.41%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.002'
$ws.Range("E24").Value = '  -0.20%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.892'
$ws.Range("E25").Value = '  +0.59%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '142.12'
$ws.Range("E26").Value = '  +0.75%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1255'
$ws.Range("E27").Value = '  +4.94%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.59'
$ws.Range("E28").Value = '  +0.32%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.756'
$ws.Range("E29").Value = '  -1.20%  '
# Row 30
$ws.Range("E30").Value = '  +0.40%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04892'
$ws.Range("E31").Value = '  +0.36%  '
# Row 32
$ws.Range("E32").Value = '  -0.80%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.181'
$ws.Range("E33").Value = '  +0.36%  '
# Row 34
$ws.Range("E34").Value = '  +0.93%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.372'
$ws.Range("E35").Value = '  +0.43%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.8968'
$ws.Range("E36").Value = '  +0.34%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.538'
$ws.Range("E37").Value = '  -1.68%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5499'
$ws.Range("E38").Value = '  +0.86%  '
# Row 39
$ws.Range("D39").Value = '1.114.30'
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01556'
$ws.Range("E40").Value = '  -0.10%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.001'
$ws.Range("E41").Value = '  +0.02%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.571'
$ws.Range("E42").Value = '  +2.56%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7951'
$ws.Range("E43").Value = '  -1.78%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '97.55'
$ws.Range("E44").Value = '  -1.58%  '
# Row 45
$ws.Range("D45").Value = '1.763.76'
$ws.Range("E45").Value = '  -0.77%  '
# Row 46
$ws.Range("D46").Value = '0.0₈117'
$ws.Range("E46").Value = '  -7.91%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4437'
$ws.Range("E47").Value = '  -2.02%  '
# Row 48
$ws.Range("E48").Value = '  +0.06%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '54.64'
# Row 50
$ws.Range("E50").Value = '  +1.46%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.527'
$ws.Range("E51").Value = '  +2.27%  '
